# Rename Description -> Product in Expenses/Revenue export sheets,
# fix Type/Item Type header colors on the Products sheet, and append
# the new expense products/services to the Products sheet.

$wb = $excel.ActiveWorkbook

# --- 1) Revenue sheet: rename D1 header "Description" -> "Product" ---
$wsRevenue = $wb.Worksheets.Item("Revenue")
$wsRevenue.Range("D1").Value = "Product"

# --- 2) Expenses sheet: rename D1 header "Description" -> "Product" ---
$wsExpenses = $wb.Worksheets.Item("Expenses")
$wsExpenses.Range("D1").Value = "Product"

# --- 3) Products sheet: fix header colors + add new rows ---
$wsProducts = $wb.Worksheets.Item("Products")

# Fix C1/D1 (Type / Item Type) header fill color to match the other
# header cells (copy format from A1, which already has the correct style).
$wsProducts.Range("A1").Copy()
$wsProducts.Range("C1:D1").PasteSpecial(-4122)

# Add the new expense products/services as rows 27-31.
$newRows = @(
    @("PRD-026", "Office Supplies Bundle", "Expenses", "Product", "OSB-001", "General office supplies", "CAT-PUR-001", "SUP-001"),
    @("PRD-027", "Printer Paper (Case)", "Expenses", "Product", "PP-500", "500 sheets per ream, 10 reams", "CAT-PUR-001", "SUP-001"),
    @("PRD-028", "Cleaning Supplies Kit", "Expenses", "Product", "CSK-001", "Commercial cleaning supplies", "CAT-PUR-002", "SUP-002"),
    @("PRD-029", "IT Support Service", "Expenses", "Service", "ITS-HR", "Hourly IT support", "CAT-PUR-003", "SUP-004"),
    @("PRD-030", "Marketing Consultation", "Expenses", "Service", "MKT-HR", "Marketing consulting hourly", "CAT-PUR-003", "SUP-003")
)

$startRow = 27
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]
    for ($j = 0; $j -lt $values.Count; $j++) {
        $col = $j + 1
        $wsProducts.Cells.Item($row, $col).Value = $values[$j]
    }
}
